# Added Edit and verification method
# - Adds createNewFlag / editFlag columns (H,I) before detailsFlag
# - Adds editCourseName column (T) at the end
# - Adds a second test case row (row 3) exercising the edit/verification flow
# - Updates the _FilterDatabase defined name + pane/selection to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (H:I) for createNewFlag / editFlag ahead of the
#    existing detailsFlag column.
# ---------------------------------------------------------------------------
$ws.Columns("H:I").Insert()

# Header row uses the bordered/filled header style (s=1) everywhere - H1/I1
# already inherit it from the insert, just set the text.
$ws.Range("H1").Value2 = "createNewFlag"
$ws.Range("I1").Value2 = "editFlag"

# ---------------------------------------------------------------------------
# 2. Row 2 clean-up: the insert leaves H2/I2 blank (inherited plain style)
#    and J2 (old H2, "detailsFlag" flag) still carrying the old one-off
#    style. Re-home the formatting from already-correct neighbour cells via
#    copy/paste-special (format only) so no stray style entries are created.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("R2").PasteSpecial(-4122)

# createNewFlag is checked for the first test case; editFlag stays blank
# (the cell is dropped entirely, matching a never-edited flag column).
$ws.Range("H2").Value2 = 1
$ws.Range("I2").Clear()

# ---------------------------------------------------------------------------
# 3. Build row 3 (second test case) by cloning row 2's formatting, then
#    filling in the new values / clearing the cells that should stay empty.
# ---------------------------------------------------------------------------
$ws.Range("A2:T2").Copy()
$ws.Range("A3:T3").PasteSpecial(-4122)
$ws.Range("A3").Select()
$excel.CutCopyMode = 0

$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "testCase_2"
$ws.Range("C3").Value2 = "Second TC"

# Append the editCourseName column header (T1) once the row-3 test case id
# and title are in place.
$ws.Range("T1").Value2 = "editCourseName"

$ws.Range("D3").Value2 = 1
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = "Course"
$ws.Range("G3").Value2 = "Courses"
$ws.Range("H3").Clear()
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = 1
$ws.Range("K3").Value2 = "TEST SERIES"

$ws.Range("M3").Value2 = "Test Auto"
$ws.Range("T3").Value2 = "Edited Course"
$ws.Range("L3").Value2 = "Test"

$ws.Range("N3").Value2 = "OPEN"
$ws.Range("O3").Value2 = "LINEAR"
$ws.Range("P3").Value2 = "SHOW"
$ws.Range("Q3").Clear()
$ws.Range("R3").ClearContents()
$ws.Range("S3").Clear()

# L3 ("Test") uses the highlighted/selection style (same as K3), not the
# plain style row-2's equivalent column (Course_1) carried over from the copy.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value2 = "Test"

$ws.Rows("3:3").RowHeight = 28.8

# ---------------------------------------------------------------------------
# 5. Sheet-level view state + the hidden _FilterDatabase defined name need to
#    track the new A1:T column span.
# ---------------------------------------------------------------------------
$names = $wb.Names
$fd = $names.Item(1)
$fd.RefersTo = "=DemoTestCaseSheet!`$A`$1:`$P`$1"

$ws.Range("O7").Select()

Write-Output "done"
